# Bulk Upload Template fix: role and technology bugs
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Technologies sheet: remove the duplicate/incorrect ".Net" row and fix the
#    last technology row (was "Java", should be "CSS").
# ---------------------------------------------------------------------------
$wsTech = $wb.Worksheets.Item("Technologies")
$wsTech.Range("B3").ClearContents()
$wsTech.Range("B7").Value2 = "CSS"

# ---------------------------------------------------------------------------
# 2) Employees sheet: fix Mrunali Desai's EMP ID (was wrongly duplicated as
#    E0122, should be E0123 like Suraksha Nigade's row used to read before
#    the fix went in).
# ---------------------------------------------------------------------------
$wsEmp = $wb.Worksheets.Item("Employees")
$wsEmp.Range("B3").Value2 = "E0123"

# ---------------------------------------------------------------------------
# 3) Roles sheet: no data changes.
# ---------------------------------------------------------------------------
$wsRoles = $wb.Worksheets.Item("Roles")

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping to mirror the saved workbook view
# state: Employees becomes the active sheet with B3 selected.
# ---------------------------------------------------------------------------
$wsTech.Range("B9").Select() | Out-Null
$wsRoles.Range("H14").Select() | Out-Null
$wsEmp.Activate() | Out-Null
$wsEmp.Range("B3").Select() | Out-Null
